$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 74-78 were re-ordered (their F:V betting-data columns were reshuffled
# among themselves) and a brand-new match (East Fife vs Bonnyrigg Rose) was
# appended as row 79. Columns A:E for rows 74-78 are untouched.

$ws.Range("F74").Value = "Stenhousemuir"
$ws.Range("G74").Value = 2
$ws.Range("H74").Value = "East Fife"
$ws.Range("I74").Value = 1
$ws.Range("J74").Value = 1.78
$ws.Range("K74").Value = "14/12/2023 09:13"
$ws.Range("L74").Value = 1.92
$ws.Range("M74").Value = "16/12/2023 15:56"
$ws.Range("N74").Value = 3.53
$ws.Range("O74").Value = "14/12/2023 09:13"
$ws.Range("P74").Value = 3.46
$ws.Range("Q74").Value = "16/12/2023 15:58"
$ws.Range("R74").Value = 4.01
$ws.Range("S74").Value = "14/12/2023 09:13"
$ws.Range("T74").Value = 4.07
$ws.Range("U74").Value = "16/12/2023 15:58"
$ws.Range("V74").Value = "https://www.betexplorer.com/football/scotland/league-two/stenhousemuir-east-fife/Q9OcVQh7/"

$ws.Range("F75").Value = "Peterhead"
$ws.Range("G75").Value = 2
$ws.Range("H75").Value = "Clyde"
$ws.Range("I75").Value = 1
$ws.Range("J75").Value = 1.55
$ws.Range("K75").Value = "14/12/2023 09:13"
$ws.Range("L75").Value = 1.54
$ws.Range("M75").Value = "16/12/2023 15:31"
$ws.Range("N75").Value = 4.01
$ws.Range("O75").Value = "14/12/2023 09:13"
$ws.Range("P75").Value = 4.29
$ws.Range("Q75").Value = "16/12/2023 15:31"
$ws.Range("R75").Value = 5.12
$ws.Range("S75").Value = "14/12/2023 09:13"
$ws.Range("T75").Value = 5.75
$ws.Range("U75").Value = "16/12/2023 15:31"
$ws.Range("V75").Value = "https://www.betexplorer.com/football/scotland/league-two/peterhead-clyde/roZhW6w1/"

$ws.Range("F76").Value = "Stranraer"
$ws.Range("G76").Value = 3
$ws.Range("H76").Value = "Bonnyrigg Rose"
$ws.Range("I76").Value = 1
$ws.Range("J76").Value = 2.34
$ws.Range("K76").Value = "14/12/2023 09:13"
$ws.Range("L76").Value = 2.67
$ws.Range("M76").Value = "16/12/2023 15:59"
$ws.Range("N76").Value = 3.16
$ws.Range("O76").Value = "14/12/2023 09:13"
$ws.Range("P76").Value = 3.26
$ws.Range("Q76").Value = "16/12/2023 15:59"
$ws.Range("R76").Value = 2.86
$ws.Range("S76").Value = "14/12/2023 09:13"
$ws.Range("T76").Value = 2.66
$ws.Range("U76").Value = "16/12/2023 15:59"
$ws.Range("V76").Value = "https://www.betexplorer.com/football/scotland/league-two/stranraer-bonnyrigg-rose/40N1Up8D/"

$ws.Range("F77").Value = "Elgin City"
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = "Spartans"
$ws.Range("I77").Value = 4
$ws.Range("J77").Value = 3.66
$ws.Range("K77").Value = "14/12/2023 09:13"
$ws.Range("L77").Value = 4.68
$ws.Range("M77").Value = "16/12/2023 15:53"
$ws.Range("N77").Value = 3.42
$ws.Range("O77").Value = "14/12/2023 09:13"
$ws.Range("P77").Value = 3.71
$ws.Range("Q77").Value = "16/12/2023 15:57"
$ws.Range("R77").Value = 1.89
$ws.Range("S77").Value = "14/12/2023 09:13"
$ws.Range("T77").Value = 1.74
$ws.Range("U77").Value = "16/12/2023 15:52"
$ws.Range("V77").Value = "https://www.betexplorer.com/football/scotland/league-two/elgin-city-spartans/zLVpYS8l/"

$ws.Range("F78").Value = "Forfar Athletic"
$ws.Range("G78").Value = 2
$ws.Range("H78").Value = "Dumbarton"
$ws.Range("I78").Value = 4
$ws.Range("J78").Value = 3.13
$ws.Range("K78").Value = "14/12/2023 09:13"
$ws.Range("L78").Value = 3.01
$ws.Range("M78").Value = "16/12/2023 15:57"
$ws.Range("N78").Value = 3.17
$ws.Range("O78").Value = "14/12/2023 09:13"
$ws.Range("P78").Value = 2.98
$ws.Range("Q78").Value = "16/12/2023 15:57"
$ws.Range("R78").Value = 2.19
$ws.Range("S78").Value = "14/12/2023 09:13"
$ws.Range("T78").Value = 2.57
$ws.Range("U78").Value = "16/12/2023 15:57"
$ws.Range("V78").Value = "https://www.betexplorer.com/football/scotland/league-two/forfar-athletic-dumbarton/UsVlXnOf/"
# Add new row 79
$ws.Range("A78").Copy($ws.Range("A79"))
$ws.Range("E78").Copy($ws.Range("E79"))
$ws.Range("A79").Value = 78
$ws.Range("B79").Value = "scotland"
$ws.Range("C79").Value = "league-two"
$ws.Range("D79").Value = "2023-2024"
$ws.Range("E79").Value = 45282.85416666666
$ws.Range("F79").Value = "East Fife"
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = "Bonnyrigg Rose"
$ws.Range("I79").Value = 3
$ws.Range("J79").Value = 1.98
$ws.Range("K79").Value = "20/12/2023 20:42"
$ws.Range("L79").Value = 2.2
$ws.Range("M79").Value = "22/12/2023 20:21"
$ws.Range("N79").Value = 3.32
$ws.Range("O79").Value = "20/12/2023 20:42"
$ws.Range("P79").Value = 3.18
$ws.Range("Q79").Value = "22/12/2023 20:20"
$ws.Range("R79").Value = 3.49
$ws.Range("S79").Value = "20/12/2023 20:42"
$ws.Range("T79").Value = 3.47
$ws.Range("U79").Value = "22/12/2023 20:21"
$ws.Range("V79").Value = "https://www.betexplorer.com/football/scotland/league-two/east-fife-bonnyrigg-rose/8WYCy2Ft/"